# Saving continuous code before making changes to plots - compare MI instead.
# Insert two new continuous-monitoring stations (Fishermans Cut / Yolo Bypass
# at Lisbon) into the stations table, just before the existing MDM row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 8-9, pushing the existing data (starting with
# the MDM station) down by two rows.
$ws.Range("A8:A9").EntireRow.Insert()

# Fill in the station codes first so the shared-string table is built in the
# same order as the source edit (StationCode column before the rest).
$ws.Range("A8").Value = "FCT"
$ws.Range("A9").Value = "LIS"

# Row 8: Fishermans Cut (FCT)
$ws.Range("C8").Value = 38.067822
$ws.Range("D8").Value = -121.648838
$ws.Range("E8").Value = "WQ"
$ws.Range("F8").Value = "FISHERMANS CUT"

# Row 9: Yolo Bypass at Lisbon (LIS)
$ws.Range("C9").Value = 38.474781
$ws.Range("D9").Value = -121.588226
$ws.Range("E9").Value = "WQ"
$ws.Range("F9").Value = "YOLO BYPASS AT LISBON"

# Match the author's final cell selection.
$ws.Range("H8").Select()
